$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.795.53'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  +0.65%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.644.89'
$ws.Range("D3").Style = "Normal"

$ws.Range("E4").Value = '  +0.45%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.98'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  +0.53%  '

$ws.Range("E7").Value = '  +0.49%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.252'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = '  -0.18%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0630'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  +0.55%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.16'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  -0.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0844'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  +0.08%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.868.69'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  -0.22%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.658.40'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  +0.29%  '

$ws.Range("E14").Value = '  -1.07%  '

$ws.Range("E15").Value = '  -0.76%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.53'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  -2.17%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.796.48'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  +0.48%  '

$ws.Range("E18").Value = '  -1.68%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '213.77'
$ws.Range("D19").Style = "Normal"

$ws.Range("E20").Value = '  +0.47%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.38'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  +0.04%  '

$ws.Range("E22").Value = '  +14.46%  '

$ws.Range("E23").Value = '  -0.74%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.36'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  -2.13%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.37'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  -0.61%  '

$ws.Range("E26").Value = '  +0.14%  '

$ws.Range("E27").Value = '  -1.30%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.67'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  -1.19%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0510'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  -1.65%  '

$ws.Range("E31").Value = '  +0.43%  '

$ws.Range("E32").Value = '  -2.10%  '

$ws.Range("E33").Value = '  -1.85%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.296.40'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  +1.72%  '

$ws.Range("E35").Value = '  -0.19%  '

$ws.Range("E36").Value = '  +1.49%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0175'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  -4.41%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.535'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  +1.00%  '

$ws.Range("E39").Value = '  -0.30%  '

$ws.Range("E40").Value = '  +0.54%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.809'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  +0.12%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.36'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  -2.01%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.795.33'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  +0.63%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.78'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  +3.23%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.63'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  -1.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.62'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  +0.85%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0526'
$ws.Range("D48").Style = "Normal"

$ws.Range("E49").Value = '  -1.32%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0977'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  +0.07%  '

$ws.Range("B51").Value = 'USDD'

$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.01'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  +0.70%  '
